$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "sports_club_coed"
$ws.Range("C2").Value = "Lacrosse"

$ws.Range("B3").Value = "sports_club_boys"
$ws.Range("C3").Value = "Lacrosse"

$ws.Range("B4").Value = "sports_club_girls"
$ws.Range("C4").Value = "Lacrosse"

$ws.Range("B5").Value = "sports_club_boys"
$ws.Range("C5").Value = "Lacrosse"

$ws.Range("B6").Value = "sports_club_boys"
$ws.Range("C6").Value = "Lacrosse"

$ws.Range("B7").Value = "sports_club_girls"
$ws.Range("C7").Value = "Lacrosse"

$ws.Range("B8").Value = "sports_club_boys"
$ws.Range("C8").Value = "Lacrosse"

$ws.Range("B9").Value = "sports_club_girls"
$ws.Range("C9").Value = "Lacrosse"

$ws.Range("B10").Value = "sports_club_coed"
$ws.Range("C10").Value = "Lacrosse"

$ws.Range("B11").Value = "sports_club_boys"
$ws.Range("C11").Value = "Lacrosse"

$ws.Range("B12").Value = "sports_club_girls"
$ws.Range("C12").Value = "Lacrosse"
